$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 62936.89
$ws.Range("J40").Value = 76905.336
$ws.Range("L40").Value = 76905.336
$ws.Range("N40").Value = -77255.336

$ws.Range("H64").Value = 2722031
$ws.Range("I64").Value = 4352275
$ws.Range("J64").Value = 4957.8335
$ws.Range("K64").Value = 4352275
$ws.Range("L64").Value = 4957.8335
$ws.Range("M64").Value = -4352027
$ws.Range("N64").Value = -5453.8335

$ws.Range("H67").Value = 2722031
$ws.Range("I67").Value = 4352275
$ws.Range("J67").Value = 4957.8335
$ws.Range("K67").Value = 4352275
$ws.Range("L67").Value = 4957.8335
$ws.Range("M67").Value = -4351417
$ws.Range("N67").Value = -6673.8335

$ws.Range("H70").Value = 4995.067
$ws.Range("J70").Value = 5954.304
$ws.Range("L70").Value = 17862.912
$ws.Range("N70").Value = -18402.912

$ws.Range("H73").Value = 4995.067
$ws.Range("J73").Value = 5954.304
$ws.Range("L73").Value = 17862.912
$ws.Range("N73").Value = -19734.912

$ws.Range("H96").Value = 2756.3635
$ws.Range("I96").Value = 538.2
$ws.Range("K96").Value = 1614.6
$ws.Range("M96").Value = -241.6000000000001

$ws.Range("H98").Value = 1089.1428
$ws.Range("I98").Value = 1256.9166
$ws.Range("J98").Value = 82.5
$ws.Range("K98").Value = 1256.9166
$ws.Range("L98").Value = 82.5
$ws.Range("M98").Value = 241.0834
$ws.Range("N98").Value = -3078.5

$ws.Range("H122").Value = 1089.1428
$ws.Range("I122").Value = 1256.9166
$ws.Range("J122").Value = 82.5
$ws.Range("K122").Value = 3770.7498
$ws.Range("L122").Value = 247.5
$ws.Range("M122").Value = -1320.7498
$ws.Range("N122").Value = -5147.5

$ws.Range("H132").Value = 10218.516
$ws.Range("I132").Value = 1564.9546
$ws.Range("K132").Value = 4694.8638
$ws.Range("M132").Value = -2164.8638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6021.186
$ws.Range("I32").Value = 3982
$ws.Range("J32").Value = 10245.214
$ws.Range("K32").Value = 3982
$ws.Range("L32").Value = 10245.214
$ws.Range("M32").Value = -3695
$ws.Range("N32").Value = -10819.214

$ws.Range("H45").Value = 3000
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2623

$ws.Range("H61").Value = 7606.2285
$ws.Range("I61").Value = 6171.8213
$ws.Range("K61").Value = 6171.8213
$ws.Range("M61").Value = -5959.8213

$ws.Range("H74").Value = 1025.303
$ws.Range("I74").Value = 649.2105
$ws.Range("J74").Value = 1535.7142
$ws.Range("K74").Value = 649.2105
$ws.Range("L74").Value = 1535.7142
$ws.Range("M74").Value = 224.7895
$ws.Range("N74").Value = -3283.7142

$ws.Range("H77").Value = 1025.303
$ws.Range("I77").Value = 649.2105
$ws.Range("J77").Value = 1535.7142
$ws.Range("K77").Value = 3246.0525
$ws.Range("L77").Value = 7678.571
$ws.Range("M77").Value = 1121.9475
$ws.Range("N77").Value = -16414.571

$ws.Range("H122").Value = 4967.0645
$ws.Range("I122").Value = 2293.7778
$ws.Range("J122").Value = 8668.538
$ws.Range("K122").Value = 6881.3334
$ws.Range("L122").Value = 26005.614
$ws.Range("M122").Value = -4431.3334
$ws.Range("N122").Value = -30905.614

$ws.Range("H132").Value = 16974.863
$ws.Range("I132").Value = 21826.25
$ws.Range("K132").Value = 65478.75
$ws.Range("M132").Value = -62948.75

$ws.Range("H136").Value = 7606.2285
$ws.Range("I136").Value = 6171.8213
$ws.Range("K136").Value = 18515.4639
$ws.Range("M136").Value = -15965.4639

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 79999.5
$ws.Range("J135").Value = 79999.5
$ws.Range("L135").Value = 79999.5
$ws.Range("N135").Value = -90139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5029.9287
$ws.Range("I122").Value = 2772.4285
$ws.Range("J122").Value = 7287.4287
$ws.Range("K122").Value = 8317.2855
$ws.Range("L122").Value = 21862.2861
$ws.Range("M122").Value = -5867.2855
$ws.Range("N122").Value = -26762.2861

$ws.Range("H134").Value = 2502.9487
$ws.Range("I134").Value = 2477.6287
$ws.Range("K134").Value = 7432.886100000001
$ws.Range("M134").Value = -4897.886100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1147530.1
$ws.Range("J34").Value = 2933.3333
$ws.Range("L34").Value = 8799.999899999999
$ws.Range("N34").Value = -8967.999899999999

$ws.Range("H39").Value = 4883.8
$ws.Range("J39").Value = 4868.5
$ws.Range("L39").Value = 14605.5
$ws.Range("N39").Value = -15193.5

$ws.Range("H55").Value = 9000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H113").Value = 927.6
$ws.Range("I113").Value = 762.6667
$ws.Range("J113").Value = 1175
$ws.Range("K113").Value = 2288.0001
$ws.Range("L113").Value = 3525
$ws.Range("M113").Value = -118.0001000000002
$ws.Range("N113").Value = -7865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 777
$ws.Range("I97").Value = 808.8889
$ws.Range("J97").Value = 490
$ws.Range("K97").Value = 808.8889
$ws.Range("L97").Value = 490
$ws.Range("M97").Value = -312.8889
$ws.Range("N97").Value = -1482

$ws.Range("H102").Value = 9616.941000000001
$ws.Range("I102").Value = 9534.857
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 9534.857
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -7912.857
$ws.Range("N102").Value = -13244

$ws.Range("H122").Value = 848846.5600000001
$ws.Range("I122").Value = 1225278.5
$ws.Range("J122").Value = 1874.75
$ws.Range("K122").Value = 3675835.5
$ws.Range("L122").Value = 5624.25
$ws.Range("M122").Value = -3673385.5
$ws.Range("N122").Value = -10524.25

$ws.Range("H126").Value = 3391.476
$ws.Range("I126").Value = 2686.1875
$ws.Range("K126").Value = 8058.5625
$ws.Range("M126").Value = -5588.5625

$ws.Range("H132").Value = 3156.6428
$ws.Range("I132").Value = 2672.0908
$ws.Range("J132").Value = 4933.3335
$ws.Range("K132").Value = 8016.2724
$ws.Range("L132").Value = 14800.0005
$ws.Range("M132").Value = -5486.2724
$ws.Range("N132").Value = -19860.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7102.1113
$ws.Range("J40").Value = 7382.2
$ws.Range("L40").Value = 7382.2
$ws.Range("N40").Value = -7654.2

$ws.Range("H93").Value = 1226.7142
$ws.Range("J93").Value = 1673.5
$ws.Range("L93").Value = 1673.5
$ws.Range("N93").Value = -4169.5

$ws.Range("H122").Value = 35718784
$ws.Range("I122").Value = 4999.5
$ws.Range("J122").Value = 71432570
$ws.Range("K122").Value = 14998.5
$ws.Range("L122").Value = 214297710
$ws.Range("M122").Value = -12548.5
$ws.Range("N122").Value = -214302610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7274.5557
$ws.Range("I62").Value = 8235.5
$ws.Range("K62").Value = 8235.5
$ws.Range("M62").Value = -7611.5

$ws.Range("H65").Value = 7274.5557
$ws.Range("I65").Value = 8235.5
$ws.Range("K65").Value = 41177.5
$ws.Range("M65").Value = -38057.5

$ws.Range("H122").Value = 4288.8374
$ws.Range("I122").Value = 3544.1333
$ws.Range("J122").Value = 6007.385
$ws.Range("K122").Value = 10632.3999
$ws.Range("L122").Value = 18022.155
$ws.Range("M122").Value = -8182.3999
$ws.Range("N122").Value = -22922.155
